$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 13 ("5840560 - Marco Antonio Carvalho Pereira" row), shifting rows 14-24 up to 13-23.
$ws.Rows("13:13").Delete()

# After the shift, update the content of several cells to the new (post-edit) values.

# Row 10 (was "Objetivos:" row) now carries the professor name instead of the long objectives text.
$ws.Range("B10").Value = "5840560 - Marco Antonio Carvalho Pereira"
$ws.Range("C10").Value = "5840560 - Marco Antonio Carvalho Pereira"

# Row 13 (was "Programa resumido:" row) now holds "Semestral" instead of the old short-syllabus text.
$ws.Range("B13").Value = "Semestral"
$ws.Range("C13").Value = "Semestral"

# Row 15 (was "Programa:" row) now holds a date instead of the old program text.
# Copy from the existing "01/01/2015" cells (row 8) so Excel keeps it as text
# (shared string) instead of auto-converting the literal into a date serial number.
$ws.Range("B8").Copy($ws.Range("B15"))
$ws.Range("C8").Copy($ws.Range("C15"))

# Row 18 (was "Método:" row) now holds the professor name.
$ws.Range("B18").Value = "5840560 - Marco Antonio Carvalho Pereira"
$ws.Range("C18").Value = "5840560 - Marco Antonio Carvalho Pereira"

# Row 19 (was "Critério:" row) now holds the old "Método:" description text.
$metodoText = @"
O método utilizado tem por fundamento a Aprendizagem Baseada em Projetos (PBL) que visa desenvolver as competências técnicas relativas ao tema do projeto, bem como competências transversais, tais como: aprender a aprender, trabalho em equipe, relacionamento interpessoal, aspectos de liderança e capacidade de comunicação, dentre outras.

Os alunos serão divididos em grupos que desenvolverão um projeto durante o semestre relacionado a um tema de Engenharia de Produção, similar ao que eles irão encontrar na vida real no efetivo exercício de sua profissão. 
Cada grupo deverá buscar o conhecimento prático necessário para ser aplicado no desenvolvimento do projeto.
As aulas ocorrerão: 1) através de uma reunião da equipe de trabalho para tratar do projeto, e  2) palestras e dinâmicas relativas ao tema do projeto, conduzidas por professores  ou profissionais de empresas.
"@
$ws.Range("B19").Value = $metodoText
$ws.Range("C19").Value = $metodoText

# Row 20 (was "Norma de recuperação:" row) now holds the old "Critério:" description text.
$criterioText = @"
A nota será individual e será a média ponderada de componentes do projeto, tais como: Projeto Preliminar, Projeto Final, envolvimento do aluno com o projeto, Avaliação dos Pares, Apresentação de Trabalhos, dentre outros.
O detalhamento dos pesos para ponderação da média da disciplina será definido por uma equipe de professores que atuarão na coordenação da disciplina.
"@
$ws.Range("B20").Value = $criterioText
$ws.Range("C20").Value = $criterioText

# Row 21 (was "Bibliografia:" row) now holds the old "Não há recuperação" text.
$ws.Range("B21").Value = "Não há recuperação"
$ws.Range("C21").Value = "Não há recuperação"
